# question_import_template.xlsx: fix import question and has Role
# - fix casing/typo on "Danh sách nhóm câu hỏi" sheet (was "Danh sách Nhóm câu hỏi" / "Mã ")
# - tidy header cell formatting on that sheet (drop stray number format / border)
# - add a new "Mức độ câu hỏi" sheet (cloned from the fixed nhóm câu hỏi sheet)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Fix content + formatting on the existing "Danh sách nhóm câu hỏi" sheet
# ---------------------------------------------------------------------------
$wsGroup = $wb.Worksheets.Item(2)

# Title text casing fix
$wsGroup.Range("A1").Value = "Danh sách nhóm câu hỏi"

# Header text: "Mã " -> "Mã" (trailing space removed)
$wsGroup.Range("B2").Value = "Mã"

# "Mã môn học" header had a stray numeric format + extra border - make it General
# and drop the number format so it behaves like the other text headers.
$wsGroup.Range("D2").NumberFormat = "General"

# Drop the right-hand border edge on the "Mã" column cells (B2/B3) so the header
# band + body band no longer double up a border against column C.
$wsGroup.Range("B2").Borders.Item(10).LineStyle = -4142
$wsGroup.Range("B3").Borders.Item(10).LineStyle = -4142

# The subject-code formula cell (D3) no longer needs right alignment / number format
$wsGroup.Range("D3").NumberFormat = "General"
$wsGroup.Range("D3").HorizontalAlignment = -4131

# cosmetic: selection on this sheet moved
$wsGroup.Range("F31").Select()

# ---------------------------------------------------------------------------
# 2) Clone this sheet to create the new "Mức độ câu hỏi" sheet (keeps fonts,
#    fills, borders, comments & their bold "author:" prefix formatting)
# ---------------------------------------------------------------------------
$wsGroup.Copy($null, $wsGroup)
$wsLevel = $wb.Worksheets.Item(3)
$wsLevel.Name = "Mức độ câu hỏi"

# Title text for the new sheet
$wsLevel.Range("A1").Value = "Danh sách mức độ câu hỏi"

# This sheet only has 3 real columns (STT / Mã / Tên); shrink the title merge
# from A1:D1 down to A1:C1 and clear out the leftover 4th ("Mã môn học") column.
$wsLevel.Range("A1:D1").UnMerge()
$wsLevel.Range("A1:C1").Merge()

$wsLevel.Range("D2").ClearContents()
$wsLevel.Range("D3").ClearContents()

# Drop the fill/border "band" styling that leaked into the now-unused column D
$wsLevel.Range("D2").Borders.Item(7).LineStyle = -4142
$wsLevel.Range("D2").Borders.Item(8).LineStyle = -4142
$wsLevel.Range("D2").Borders.Item(9).LineStyle = -4142
$wsLevel.Range("D2").Borders.Item(10).LineStyle = -4142
$wsLevel.Range("D2").Interior.ColorIndex = -4142

$wsLevel.Range("D3").Borders.Item(7).LineStyle = -4142
$wsLevel.Range("D3").Borders.Item(8).LineStyle = -4142
$wsLevel.Range("D3").Borders.Item(9).LineStyle = -4142
$wsLevel.Range("D3").Borders.Item(10).LineStyle = -4142
$wsLevel.Range("D3").Interior.ColorIndex = -4142

# D1 sits outside the shrunk title merge now - remove its title-bar fill/border
# and center alignment, keep it simply vertically centered
$wsLevel.Range("D1").Interior.ColorIndex = -4142
$wsLevel.Range("D1").Borders.Item(7).LineStyle = -4142
$wsLevel.Range("D1").Borders.Item(8).LineStyle = -4142
$wsLevel.Range("D1").Borders.Item(9).LineStyle = -4142
$wsLevel.Range("D1").Borders.Item(10).LineStyle = -4142
$wsLevel.Range("D1").HorizontalAlignment = -4142
$wsLevel.Range("D1").VerticalAlignment = -4108

# A1/B1/C1 pick up the merge-edge borders now that the merge ends at C1
$wsLevel.Range("A1").Borders.Item(9).LineStyle = 1
$wsLevel.Range("B1").Borders.Item(9).LineStyle = 1
$wsLevel.Range("C1").Borders.Item(9).LineStyle = 1
$wsLevel.Range("C1").Borders.Item(10).LineStyle = 1

# Update the two legacy comments copied over from the nhóm câu hỏi sheet so they
# describe this sheet's own jxls loop (questionTypes / lastCell C3 instead of
# groupQuestions / D3), keeping the bold "author:" prefix formatting intact.
$cmt1 = $wsLevel.Range("A1").Comment
$cmt1.Text("sang nguyen:`njx:area(lastCell = ""C3"")")

$cmt3 = $wsLevel.Range("A3").Comment
$cmt3.Text("admin:`njx:each(items=""questionTypes"", var=""item"", varIndex=""i"", lastCell=""C3"")")

# cosmetic: selection/view on the new sheet
$wsLevel.Range("D24").Select()
